$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3293
$ws1.Range("F5").Value = 2378
$ws1.Range("F7").Value = 331
$ws1.Range("F8").Value = 1347
$ws1.Range("F9").Value = 1065
$ws1.Range("F10").Value = 284
$ws1.Range("F11").Value = 501
$ws1.Range("F16").Value = 8303
$ws1.Range("F23").Value = 570
$ws1.Range("F27").Value = 1922
$ws1.Range("F28").Value = 1604
$ws1.Range("F29").Value = 63
$ws1.Range("F33").Value = 17
$ws1.Range("F34").Value = 23
$ws1.Range("F38").Value = 53
$ws1.Range("F39").Value = 215
$ws1.Range("F40").Value = 390
$ws1.Range("F41").Value = 60

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3293
$ws4.Range("F6").Value = 2378
$ws4.Range("F8").Value = 331
$ws4.Range("F9").Value = 1347
$ws4.Range("F11").Value = 1065
$ws4.Range("F12").Value = 284
$ws4.Range("F13").Value = 501
$ws4.Range("F17").Value = 8303
$ws4.Range("F25").Value = 570
$ws4.Range("F29").Value = 1922
$ws4.Range("F30").Value = 1605
$ws4.Range("F34").Value = 17
$ws4.Range("F35").Value = 23
$ws4.Range("F39").Value = 53
$ws4.Range("F40").Value = 215
$ws4.Range("F41").Value = 390
$ws4.Range("F46").Value = 60
